# #1440 change surveySeries to studySeries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")

$ws.Range("H1").Value = "studySeries.en"
$ws.Range("G1").Value = "studySeries.de"

$ws.Range("E2").Select()
